# Generate Report for Handback
#
# Marks the a.md / b.md files as handed back (in sync with en-US) for
# both the zh-cn and de-de localization targets, records the handback
# xliff file name + timestamp for each language, adds a hyperlink on
# the newly-populated "Latest Target File" cell, and widens the columns
# that now hold the longer status/file-name text.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67378b97842fc3df53fc5395ed1f1d9127203966/e2e/"
$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns for zh-cn (E) and de-de (F) ---------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn detail sheet ---------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

# Re-create all hyperlinks on the sheet so the new "Latest Target File"
# link (column I) is interleaved with the existing "Source File Name"
# links (column A) in row order, matching rows 2 and 3.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $repoBase + "a.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $repoBase + "a.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $repoBase + "b.md", "", "", "b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $repoBase + "a.md", "", "", "a.md")

# Latest Handback File (column J) references the handoff xliff for zh-cn.
$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# Latest Handback DateTime (column K).
$wsZhCn.Range("K2").Value = "2016-08-17 10:34:35"
$wsZhCn.Range("K3").Value = "2016-08-17 10:34:35"

# --- de-de detail sheet ----------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $repoBase + "a.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $repoBase + "a.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $repoBase + "b.md", "", "", "b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $repoBase + "a.md", "", "", "a.md")

$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-08-17 10:34:42"
$wsDeDe.Range("K3").Value = "2016-08-17 10:34:42"

# --- Widen columns that now hold the longer text --------------------------
# (ColumnWidth is specified in character units; this runtime snaps the
# resulting value to the nearest 1/6 of a character, so we request the
# input that lands exactly on the desired output width.)
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 29.166666666666668
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 29.166666666666668

$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 29.166666666666668
$wsZhCn.Range("J1").EntireColumn.ColumnWidth = 39.166666666666664

$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 29.166666666666668
$wsDeDe.Range("J1").EntireColumn.ColumnWidth = 39.166666666666664
